$d = $word.ActiveDocument

# 1. Update "Total de citas programadas: 4" -> "...: 2"
#    The number lives in its own run (separate from the bold label run),
#    so we target just that trailing character instead of using Find
#    (which would merge the two runs into one).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Total de citas programadas:*") {
        $e = $p.Range.End
        $numRange = $d.Range($e - 2, $e - 1)
        if ($numRange.Text -eq "4") {
            $numRange.Text = "2"
        }
    }
}

# 2. Remove the last two appointment rows from the schedule table
#    (09:45 - 10:00 / REGIONAL S.A.S  and  10:15 - 10:30 / BOX BRAND),
#    leaving the header row plus the PROCOLOMBIA and INTERLINK2AMERICAS rows.
$t = $d.Tables.Item(1)
while ($t.Rows.Count -gt 3) {
    $t.Rows.Item($t.Rows.Count).Delete()
}
